$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.22420346736908
$ws.Range("B1").Value = 2.699114561080933
$ws.Range("C1").Value = 4.250983238220215
$ws.Range("D1").Value = 2.138493537902832
$ws.Range("E1").Value = 1.157812118530273
